$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Force text format on every touched cell first so Excel does not
# auto-coerce numeric-looking strings (e.g. "4.996", "0.000008680")
# into numbers/dates and strip significant trailing zeros / dot-grouping.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.111.47'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.835.16'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.87'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6598'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.77%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.70'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +6.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07376'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2943'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.08'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07717'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.835.51'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.996'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6695'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.01'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.122'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008680'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.130.11'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.091.78'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.47'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.83'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.54%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.141'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.66%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.30'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.565'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.43%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.96'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.509'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.123'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.031'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.204'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05391'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.844'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7442'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.157'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.654'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.295.86'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.764'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01793'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.347'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9016'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.44%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.990.39'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07833'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.15%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000123'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '64.75'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.47%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.745'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.50%  '
